# Change the table border color from black (000000) to gray (666666)
# on the outer horizontal rules of the "poor tracking summary" table:
#   - the header row's top AND bottom borders
#   - the first data row's top border (shared edge with the header)
#   - the last data row's bottom border (closing rule of the table)
#
# Word's Borders collection index mapping observed in this runtime:
#   -1 = Top, -2 = Left, -3 = Bottom, -4 = Right
# The Color property takes a Windows COLORREF-style integer (0xBBGGRR);
# for a neutral gray (666666) R=G=B so the byte order is irrelevant.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$grayColor = 0x666666
$wdBorderTop = -1
$wdBorderBottom = -3

$colCount = $t.Columns.Count
$rowCount = $t.Rows.Count

# Header row (row 1): both top and bottom borders turn gray.
for ($col = 1; $col -le $colCount; $col++) {
    $cell = $t.Cell(1, $col)
    $cell.Borders.Item($wdBorderTop).Color = $grayColor
    $cell.Borders.Item($wdBorderBottom).Color = $grayColor
}

# First data row (row 2): its top border (the line under the header) turns gray.
for ($col = 1; $col -le $colCount; $col++) {
    $cell = $t.Cell(2, $col)
    $cell.Borders.Item($wdBorderTop).Color = $grayColor
}

# Last data row: its bottom border (closing rule of the table) turns gray.
for ($col = 1; $col -le $colCount; $col++) {
    $cell = $t.Cell($rowCount, $col)
    $cell.Borders.Item($wdBorderBottom).Color = $grayColor
}

Write-Host "Updated borders on" $colCount "columns across header/first/last rows (of" $rowCount "rows)."
